# TC_53557 - add "DC Unit Loading Details Name" / "Current (DC Units)" header
# cells in F1:F2 (matching the existing row7/row8 header & data styling) and
# move the sheet's active selection to B7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1 gets the same look as the other row-7 column headers (bold header style)
$ws.Range("F7").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "DC Unit Loading Details Name"

# F2 gets the same look as the row-8 data cells
$ws.Range("A8").Copy() | Out-Null
$ws.Range("F2").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Value = "Current (DC Units)"

# clear marching-ants / copy mode and update the saved selection to B7
$excel.CutCopyMode = 0
$ws.Range("B7").Select() | Out-Null
